$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "without GDP: B = NA, 95% CI [NA, NA],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "without GDP: B = 51.61, 95% CI [25.13, 78.09],", 1)

$d.Content.Find.Execute(
    "= NA; with GDP: B = 51.45, 95% CI [18.56, 84.33],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= < 0.001; with GDP: B = 51.45, 95% CI [18.56, 84.33],", 1)

$d.Content.Find.Execute(
    "= 0.045) or not (B = NA, 95% CI [NA, NA],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.045) or not (B = 32.08, 95% CI [-95.10, 159.25],", 1)

$d.Content.Find.Execute(
    "= NA).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.621).", 1)

$d.Content.Find.Execute(
    "For models evaluating Positivivity Bias, a significantly positive relationship was found between Positivity Bias and Collectivistic Values (B = NA, 95% CI [NA, NA],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For models evaluating Positivivity Bias, a significantly positive relationship was found between Positivity Bias and Collectivistic Values (B = 233.35, 95% CI [0.58, 466.12],", 1)

$d.Content.Find.Execute(
    "= NA) but the significance diminished when GDP per capita was controlled for (B = 238.21, 95% CI [-6.52, 482.94],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.049) but the significance diminished when GDP per capita was controlled for (B = 238.21, 95% CI [-6.52, 482.94],", 1)

$d.Content.Find.Execute(
    "= 0.056). Cultural Tightness showed strong positive correlation with Positivity Bias before (B = NA, 95% CI [NA, NA],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.056). Cultural Tightness showed strong positive correlation with Positivity Bias before (B = 27.43, 95% CI [-55.87, 110.73],", 1)

$d.Content.Find.Execute(
    "= NA) and after controlling for GDP per capita (B = 24.75, 95% CI [-62.80, 112.30],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.519) and after controlling for GDP per capita (B = 24.75, 95% CI [-62.80, 112.30],", 1)

$d.Content.Find.Execute(
    "Focusing on Warmth Bias, both Collectivist Values and Cultural Tightness were positively correlated with Warmth Bias towards older adults before (Collectivistic: B = NA, 95% CI [NA, NA],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Focusing on Warmth Bias, both Collectivist Values and Cultural Tightness were positively correlated with Warmth Bias towards older adults before (Collectivistic: B = 339.01, 95% CI [34.26, 643.76],", 1)

$d.Content.Find.Execute(
    "= NA; Tightness: B = NA, 95% CI [NA, NA],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.029; Tightness: B = 38.03, 95% CI [-101.24, 177.30],", 1)

$d.Content.Find.Execute(
    "= NA) or after controlling for the influence of GDP per capita (Collectivistic: B = 395.58, 95% CI [6.82, 784.35],",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "= 0.593) or after controlling for the influence of GDP per capita (Collectivistic: B = 395.58, 95% CI [6.82, 784.35],", 1)
